$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) starting at D1
$ws.Range("D1").Value = "Zone"
$ws.Range("E1").Value = "Sector"
$ws.Range("F1").Value = "FAMILY"
$ws.Range("G1").Value = "COMMUNICATION"
$ws.Range("H1").Value = "WORK/FAMILY BALANCE"
$ws.Range("I1").Value = "WORK-RESPECT"
$ws.Range("J1").Value = "MORALIS ASPECT"
$ws.Range("K1").Value = "KERALA FESTIVAl "
$ws.Range("L1").Value = "HEALTH"
$ws.Range("M1").Value = "EXERCISE"
$ws.Range("N1").Value = "AIM/DREAM"
$ws.Range("O1").Value = "AS A SOCIAL WORKER"

# Copy style from an existing header cell (C1) to the new header cells (D1:O1)
$ws.Range("C1").Copy()
$ws.Range("D1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update data row (row 2) starting at D2
$ws.Range("D2").Value = "RIYADH NORTH"
$ws.Range("E2").Value = "MURSALATH"
$ws.Range("F2").Value = "Maybe"
$ws.Range("G2").Value = "Multiple Times"
$ws.Range("H2").Value = "Satisfied"
$ws.Range("I2").Value = "No"
$ws.Range("J2").Value = "No"
$ws.Range("K2").Value = "No"
$ws.Range("L2").Value = "No"
$ws.Range("M2").Value = "Yes"
$ws.Range("N2").Value = "Yes"
$ws.Range("O2").Value = "Yes"
